$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 63877.4
$ws.Range("J87").Value = 69641.55499999999
$ws.Range("L87").Value = 69641.55499999999
$ws.Range("N87").Value = -72137.55499999999
$ws.Range("H90").Value = 63877.4
$ws.Range("J90").Value = 69641.55499999999
$ws.Range("L90").Value = 208924.665
$ws.Range("N90").Value = -221404.665
$ws.Range("H92").Value = 132.25
$ws.Range("J92").Value = 225
$ws.Range("L92").Value = 225
$ws.Range("N92").Value = -2721
$ws.Range("H100").Value = 4999.5
$ws.Range("I100").Value = 3499.25
$ws.Range("K100").Value = 3499.25
$ws.Range("M100").Value = -2958.25
$ws.Range("H113").Value = 1446.4
$ws.Range("I113").Value = 1446.2142
$ws.Range("J113").Value = 1449
$ws.Range("K113").Value = 1446.2142
$ws.Range("L113").Value = 1449
$ws.Range("M113").Value = 1807.7858
$ws.Range("N113").Value = -7957
$ws.Range("H127").Value = 4943.5
$ws.Range("I127").Value = 4935.4287
$ws.Range("K127").Value = 14806.2861
$ws.Range("M127").Value = -9846.286100000001
$ws.Range("H132").Value = 9646.607
$ws.Range("I132").Value = 8850.23
$ws.Range("K132").Value = 26550.69
$ws.Range("M132").Value = -24020.69
$ws.Range("H135").Value = 739.8333
$ws.Range("I135").Value = 221.2
$ws.Range("K135").Value = 1990.8
$ws.Range("M135").Value = 544.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 831.25
$ws.Range("I12").Value = 700
$ws.Range("K12").Value = 700
$ws.Range("M12").Value = -527
$ws.Range("H26").Value = 631.75
$ws.Range("I26").Value = 631.75
$ws.Range("K26").Value = 631.75
$ws.Range("M26").Value = -301.75
$ws.Range("H61").Value = 2307.3635
$ws.Range("I61").Value = 2307.3635
$ws.Range("K61").Value = 2307.3635
$ws.Range("M61").Value = -2095.3635
$ws.Range("H74").Value = 1086.2941
$ws.Range("I74").Value = 966.6875
$ws.Range("K74").Value = 966.6875
$ws.Range("M74").Value = -92.6875
$ws.Range("H77").Value = 1086.2941
$ws.Range("I77").Value = 966.6875
$ws.Range("K77").Value = 4833.4375
$ws.Range("M77").Value = -465.4375
$ws.Range("H131").Value = 67500
$ws.Range("J131").Value = 100000
$ws.Range("L131").Value = 100000
$ws.Range("N131").Value = -110080
$ws.Range("H132").Value = 6163.1665
$ws.Range("I132").Value = 6163.1665
$ws.Range("K132").Value = 18489.4995
$ws.Range("M132").Value = -15959.4995
$ws.Range("H136").Value = 2307.3635
$ws.Range("I136").Value = 2307.3635
$ws.Range("K136").Value = 6922.0905
$ws.Range("M136").Value = -4372.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 11875625
$ws.Range("I7").Value = 5000200
$ws.Range("J7").Value = 23334666
$ws.Range("K7").Value = 5000200
$ws.Range("L7").Value = 23334666
$ws.Range("M7").Value = -5000087
$ws.Range("N7").Value = -23334892
$ws.Range("H134").Value = 6528
$ws.Range("I134").Value = 4635.3335
$ws.Range("J134").Value = 9367
$ws.Range("K134").Value = 13906.0005
$ws.Range("L134").Value = 28101
$ws.Range("M134").Value = -11371.0005
$ws.Range("N134").Value = -33171

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5070.1665
$ws.Range("I31").Value = 2994
$ws.Range("J31").Value = 8530.444
$ws.Range("K31").Value = 2994
$ws.Range("L31").Value = 8530.444
$ws.Range("M31").Value = -2699
$ws.Range("N31").Value = -9120.444
$ws.Range("H34").Value = 5070.1665
$ws.Range("I34").Value = 2994
$ws.Range("J34").Value = 8530.444
$ws.Range("K34").Value = 2994
$ws.Range("L34").Value = 8530.444
$ws.Range("M34").Value = -2792
$ws.Range("N34").Value = -8934.444
$ws.Range("H58").Value = 2907.5715
$ws.Range("I58").Value = 2396.56
$ws.Range("K58").Value = 2396.56
$ws.Range("M58").Value = -2193.56
$ws.Range("H134").Value = 2533.625
$ws.Range("I134").Value = 1957.1072
$ws.Range("K134").Value = 5871.321599999999
$ws.Range("M134").Value = -3336.321599999999
$ws.Range("H136").Value = 2907.5715
$ws.Range("I136").Value = 2396.56
$ws.Range("K136").Value = 7189.68
$ws.Range("M136").Value = -4639.68

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2716.3333
$ws.Range("J34").Value = 5333.3335
$ws.Range("L34").Value = 16000.0005
$ws.Range("N34").Value = -16168.0005
$ws.Range("H39").Value = 7551.375
$ws.Range("J39").Value = 9999.5
$ws.Range("L39").Value = 29998.5
$ws.Range("N39").Value = -30586.5
$ws.Range("H55").Value = 2850
$ws.Range("J55").Value = 4875
$ws.Range("L55").Value = 14625
$ws.Range("N55").Value = -14979
$ws.Range("H68").Value = 1000
$ws.Range("J68").Value = 1033.3334
$ws.Range("L68").Value = 3100.0002
$ws.Range("N68").Value = -4722.0002
$ws.Range("H71").Value = 1000
$ws.Range("J71").Value = 1033.3334
$ws.Range("L71").Value = 9300.000599999999
$ws.Range("N71").Value = -17412.0006
$ws.Range("H136").Value = 5228.5
$ws.Range("I136").Value = 459.5
$ws.Range("J136").Value = 9997.5
$ws.Range("K136").Value = 1378.5
$ws.Range("L136").Value = 29992.5
$ws.Range("M136").Value = 3721.5
$ws.Range("N136").Value = -40192.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5428.7144
$ws.Range("I80").Value = 3600.4
$ws.Range("J80").Value = 9999.5
$ws.Range("K80").Value = 3600.4
$ws.Range("L80").Value = 9999.5
$ws.Range("M80").Value = -2602.4
$ws.Range("N80").Value = -11995.5
$ws.Range("H83").Value = 5428.7144
$ws.Range("I83").Value = 3600.4
$ws.Range("J83").Value = 9999.5
$ws.Range("K83").Value = 18002
$ws.Range("L83").Value = 49997.5
$ws.Range("M83").Value = -13010
$ws.Range("N83").Value = -59981.5
$ws.Range("H113").Value = 7661.875
$ws.Range("I113").Value = 6086.25
$ws.Range("K113").Value = 6086.25
$ws.Range("M113").Value = -3916.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8874.25
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 8874.25
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 8874.25
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -9098.25
$ws.Range("H100").Value = 5737.913
$ws.Range("I100").Value = 2270.182
$ws.Range("K100").Value = 2270.182
$ws.Range("M100").Value = -1729.182
$ws.Range("H126").Value = 8874.25
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 8874.25
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 26622.75
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -31562.75
$ws.Range("H132").Value = 3849.625
$ws.Range("J132").Value = 9999.5
$ws.Range("L132").Value = 29998.5
$ws.Range("N132").Value = -35058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H81").Value = 1781.1428
$ws.Range("I81").Value = 1781.1428
$ws.Range("K81").Value = 3562.2856
$ws.Range("M81").Value = -2501.2856
$ws.Range("H84").Value = 1781.1428
$ws.Range("I84").Value = 1781.1428
$ws.Range("K84").Value = 17811.428
$ws.Range("M84").Value = -12507.428
$ws.Range("H107").Value = 789.44446
$ws.Range("I107").Value = 734.375
$ws.Range("J107").Value = 1230
$ws.Range("K107").Value = 2203.125
$ws.Range("L107").Value = 3690
$ws.Range("M107").Value = -283.125
$ws.Range("N107").Value = -7530
$ws.Range("H126").Value = 3481.682
$ws.Range("I126").Value = 2168.5625
$ws.Range("K126").Value = 6505.6875
$ws.Range("M126").Value = -4035.6875
$ws.Range("H132").Value = 3597.3142
$ws.Range("I132").Value = 3454.1538
$ws.Range("K132").Value = 10362.4614
$ws.Range("M132").Value = -7832.4614
$ws.Range("H136").Value = 4030.2632
$ws.Range("I136").Value = 2527.6
$ws.Range("J136").Value = 6920
$ws.Range("K136").Value = 7582.799999999999
$ws.Range("L136").Value = 20760
$ws.Range("M136").Value = -5032.799999999999
$ws.Range("N136").Value = -25860

Write-Output "Applied 208 cell updates across 8 sheets"